# Updates Gilgamesh_Profits market data cells per scheduled runner refresh
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4582
$ws.Range("J40").Value = 4752
$ws.Range("L40").Value = 4752
$ws.Range("N40").Value = -5102
$ws.Range("H51").Value = 116669416
$ws.Range("I51").Value = 125003250
$ws.Range("J51").Value = 100001750
$ws.Range("K51").Value = 125003250
$ws.Range("L51").Value = 100001750
$ws.Range("M51").Value = -125002766
$ws.Range("N51").Value = -100002718
$ws.Range("H97").Value = 8537
$ws.Range("J97").Value = 8537
$ws.Range("L97").Value = 25611
$ws.Range("N97").Value = -26603
$ws.Range("H113").Value = 1499
$ws.Range("I113").Value = 1499
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1499
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 1755
$ws.Range("N113").ClearContents()
$ws.Range("H121").Value = 2032.3334
$ws.Range("J121").Value = 2032.3334
$ws.Range("L121").Value = 6097.0002
$ws.Range("N121").Value = -9591.0002
$ws.Range("H137").Value = 5708.52
$ws.Range("I137").Value = 1734.5264
$ws.Range("J137").Value = 18292.834
$ws.Range("K137").Value = 5203.5792
$ws.Range("L137").Value = 54878.50199999999
$ws.Range("M137").Value = -2653.5792
$ws.Range("N137").Value = -59978.50199999999
$ws.Range("H141").Value = 2985.9546
$ws.Range("J141").Value = 4726.25
$ws.Range("L141").Value = 14178.75
$ws.Range("N141").Value = -24538.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7077.9766
$ws.Range("I32").Value = 5712.4326
$ws.Range("J32").Value = 15498.833
$ws.Range("K32").Value = 5712.4326
$ws.Range("L32").Value = 15498.833
$ws.Range("M32").Value = -5425.4326
$ws.Range("N32").Value = -16072.833
$ws.Range("H61").Value = 3770
$ws.Range("I61").Value = 2072.8
$ws.Range("J61").Value = 6855.8184
$ws.Range("K61").Value = 2072.8
$ws.Range("L61").Value = 6855.8184
$ws.Range("M61").Value = -1860.8
$ws.Range("N61").Value = -7279.8184
$ws.Range("H74").Value = 163195.4
$ws.Range("I74").Value = 242474.56
$ws.Range("K74").Value = 242474.56
$ws.Range("M74").Value = -241600.56
$ws.Range("H77").Value = 163195.4
$ws.Range("I77").Value = 242474.56
$ws.Range("K77").Value = 1212372.8
$ws.Range("M77").Value = -1208004.8
$ws.Range("H122").Value = 3604.3076
$ws.Range("I122").Value = 3370.9768
$ws.Range("J122").Value = 4719.1113
$ws.Range("K122").Value = 10112.9304
$ws.Range("L122").Value = 14157.3339
$ws.Range("M122").Value = -7662.930399999999
$ws.Range("N122").Value = -19057.3339
$ws.Range("H136").Value = 3770
$ws.Range("I136").Value = 2072.8
$ws.Range("J136").Value = 6855.8184
$ws.Range("K136").Value = 6218.400000000001
$ws.Range("L136").Value = 20567.4552
$ws.Range("M136").Value = -3668.400000000001
$ws.Range("N136").Value = -25667.4552

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4065.24
$ws.Range("I86").Value = 3884.389
$ws.Range("K86").Value = 3884.389
$ws.Range("M86").Value = -2761.389
$ws.Range("H89").Value = 4065.24
$ws.Range("I89").Value = 3884.389
$ws.Range("K89").Value = 19421.945
$ws.Range("M89").Value = -13805.945
$ws.Range("H134").Value = 3029.9666
$ws.Range("I134").Value = 2321.0435
$ws.Range("K134").Value = 6963.130500000001
$ws.Range("M134").Value = -4428.130500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 66668024
$ws.Range("I7").Value = 1538.4
$ws.Range("K7").Value = 1538.4
$ws.Range("M7").Value = -1425.4
$ws.Range("H22").Value = 1230.7778
$ws.Range("I22").Value = 953.0909
$ws.Range("K22").Value = 953.0909
$ws.Range("M22").Value = -603.0909
$ws.Range("H31").Value = 3690.8462
$ws.Range("I31").Value = 3065.2068
$ws.Range("J31").Value = 4194.8335
$ws.Range("K31").Value = 3065.2068
$ws.Range("L31").Value = 4194.8335
$ws.Range("M31").Value = -2770.2068
$ws.Range("N31").Value = -4784.8335
$ws.Range("H34").Value = 3690.8462
$ws.Range("I34").Value = 3065.2068
$ws.Range("J34").Value = 4194.8335
$ws.Range("K34").Value = 3065.2068
$ws.Range("L34").Value = 4194.8335
$ws.Range("M34").Value = -2863.2068
$ws.Range("N34").Value = -4598.8335
$ws.Range("H56").Value = 15541.5
$ws.Range("J56").Value = 20000
$ws.Range("L56").Value = 20000
$ws.Range("N56").Value = -21690
$ws.Range("H106").Value = 369499.75
$ws.Range("J106").Value = 369499.75
$ws.Range("L106").Value = 369499.75
$ws.Range("N106").Value = -372023.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 3933
$ws.Range("I63").Value = 1800
$ws.Range("J63").Value = 4999.5
$ws.Range("K63").Value = 5400
$ws.Range("L63").Value = 14998.5
$ws.Range("M63").Value = -4651
$ws.Range("N63").Value = -16496.5
$ws.Range("H66").Value = 3933
$ws.Range("I66").Value = 1800
$ws.Range("J66").Value = 4999.5
$ws.Range("K66").Value = 16200
$ws.Range("L66").Value = 44995.5
$ws.Range("M66").Value = -12456
$ws.Range("N66").Value = -52483.5
$ws.Range("H132").Value = 2352.0286
$ws.Range("J132").Value = 2874.8096
$ws.Range("L132").Value = 25873.2864
$ws.Range("N132").Value = -30933.2864

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 29135
$ws.Range("I58").Value = 29346.666
$ws.Range("J58").Value = 28500
$ws.Range("K58").Value = 29346.666
$ws.Range("L58").Value = 28500
$ws.Range("M58").Value = -29069.666
$ws.Range("N58").Value = -29054
$ws.Range("H70").Value = 15692293
$ws.Range("I70").Value = 22822562
$ws.Range("K70").Value = 22822562
$ws.Range("M70").Value = -22822292
$ws.Range("H73").Value = 15692293
$ws.Range("I73").Value = 22822562
$ws.Range("K73").Value = 22822562
$ws.Range("M73").Value = -22821626
$ws.Range("H80").Value = 40001670
$ws.Range("I80").Value = 83334730
$ws.Range("K80").Value = 83334730
$ws.Range("M80").Value = -83333732
$ws.Range("H83").Value = 40001670
$ws.Range("I83").Value = 83334730
$ws.Range("K83").Value = 416673650
$ws.Range("M83").Value = -416668658
$ws.Range("H106").Value = 49000
$ws.Range("J106").Value = 49000
$ws.Range("L106").Value = 49000
$ws.Range("N106").Value = -51524
$ws.Range("H132").Value = 2324.976
$ws.Range("I132").Value = 2312.8235
$ws.Range("J132").Value = 2376.625
$ws.Range("K132").Value = 6938.470499999999
$ws.Range("L132").Value = 7129.875
$ws.Range("M132").Value = -4408.470499999999
$ws.Range("N132").Value = -12189.875
$ws.Range("H136").Value = 10617.223
$ws.Range("J136").Value = 10617.223
$ws.Range("L136").Value = 31851.669
$ws.Range("N136").Value = -36951.669
$ws.Range("H139").Value = 96125.55499999999
$ws.Range("J139").Value = 96125.55499999999
$ws.Range("L139").Value = 96125.55499999999
$ws.Range("N139").Value = -106405.555

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H9").Value = 5400
$ws.Range("I9").Value = 200
$ws.Range("J9").Value = 8000
$ws.Range("K9").Value = 200
$ws.Range("L9").Value = 8000
$ws.Range("M9").Value = 24
$ws.Range("N9").Value = -8448
$ws.Range("H22").Value = 596
$ws.Range("I22").Value = 595.4
$ws.Range("J22").Value = 596.75
$ws.Range("K22").Value = 595.4
$ws.Range("L22").Value = 596.75
$ws.Range("M22").Value = -300.4
$ws.Range("N22").Value = -1186.75
$ws.Range("H27").Value = 596
$ws.Range("I27").Value = 595.4
$ws.Range("J27").Value = 596.75
$ws.Range("K27").Value = 595.4
$ws.Range("L27").Value = 596.75
$ws.Range("M27").Value = -488.4
$ws.Range("N27").Value = -810.75
$ws.Range("H40").Value = 5278.0557
$ws.Range("I40").Value = 5429.533
$ws.Range("K40").Value = 5429.533
$ws.Range("M40").Value = -5293.533
$ws.Range("H98").Value = 142666.67
$ws.Range("J98").Value = 142666.67
$ws.Range("L98").Value = 142666.67
$ws.Range("N98").Value = -148656.67
$ws.Range("H122").Value = 3528.2666
$ws.Range("J122").Value = 2899.2
$ws.Range("L122").Value = 8697.599999999999
$ws.Range("N122").Value = -13597.6
$ws.Range("H136").Value = 5344.8823
$ws.Range("I136").Value = 6164.1113
$ws.Range("J136").Value = 4423.25
$ws.Range("K136").Value = 18492.3339
$ws.Range("L136").Value = 13269.75
$ws.Range("M136").Value = -15942.3339
$ws.Range("N136").Value = -18369.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 830.5
$ws.Range("I107").Value = 890.7
$ws.Range("K107").Value = 2672.1
$ws.Range("M107").Value = -752.1000000000004
$ws.Range("H122").Value = 15628665
$ws.Range("I122").Value = 3625.1
$ws.Range("J122").Value = 41670400
$ws.Range("K122").Value = 10875.3
$ws.Range("L122").Value = 125011200
$ws.Range("M122").Value = -8425.299999999999
$ws.Range("N122").Value = -125016100
$ws.Range("H132").Value = 8549992
$ws.Range("I132").Value = 10104133
$ws.Range("K132").Value = 30312399
$ws.Range("M132").Value = -30309869
$ws.Range("H136").Value = 19232880
$ws.Range("I136").Value = 23810598
$ws.Range("K136").Value = 71431794
$ws.Range("M136").Value = -71429244

Write-Host "Applied all market data updates"